$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Smalltalk): Variables & Classes 1 -> 2
$ws.Range("D3").Value = 2

# Row 4 (Ruby): Variables & Classes 1 -> 2
$ws.Range("D4").Value = 2
# I4 gains the green "highlight" fill (reuse existing style used by I2/I3)
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial(-4122)  # xlPasteFormats

# Row 5 (Java): several corrections
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 2
# H5 newly populated - copy the numeric-cell formatting from G5 first, then set value
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H5").Value = 2
# I5 gains the green "highlight" fill
$ws.Range("I2").Copy()
$ws.Range("I5").PasteSpecial(-4122)  # xlPasteFormats

# Row 6 (C#): several corrections
$ws.Range("B6").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("G6").Value = 3
# H6 newly populated
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H6").Value = 2
# I6 gains the green "highlight" fill
$ws.Range("I2").Copy()
$ws.Range("I6").PasteSpecial(-4122)  # xlPasteFormats

# Row 7 (C++): corrections
$ws.Range("B7").Value = 3
$ws.Range("D7").Value = 3

# Row 8 (Python): correction
$ws.Range("D8").Value = 3

# Update active cell / selection to B7
$ws.Range("B7").Select()
